$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "1.4536340238509356e-17" -as [double]
$ws.Range("J2").Value = "6.955149978429623e-17" -as [double]
$ws.Range("K2").Value = "0.32268579387916946" -as [double]
$ws.Range("I3").Value = "-9.265182284870255e-19" -as [double]
$ws.Range("J3").Value = "1.128382069806491e-16" -as [double]
$ws.Range("K3").Value = "0.20242779362438088" -as [double]
$ws.Range("G4").Value = "9.71445146547012e-17" -as [double]
$ws.Range("H4").Value = "-1.179611963664229e-16" -as [double]
$ws.Range("I4").Value = "3.4051979343248437e-18" -as [double]
$ws.Range("J4").Value = "5.789406928771223e-17" -as [double]
$ws.Range("K4").Value = "0.26914240915726706" -as [double]
$ws.Range("I5").Value = "2.4033148156762014e-17" -as [double]
$ws.Range("J5").Value = "4.95229378484802e-17" -as [double]
$ws.Range("K5").Value = "0.27851736593571297" -as [double]
$ws.Range("H6").Value = "-6.938893903907228e-17" -as [double]
$ws.Range("I6").Value = "-4.04768811061255e-18" -as [double]
$ws.Range("J6").Value = "1.632443510633924e-17" -as [double]
$ws.Range("K6").Value = "0.2138726395983804" -as [double]
$ws.Range("I7").Value = "-2.1367665711386815e-19" -as [double]
$ws.Range("J7").Value = "6.877623732883437e-17" -as [double]
$ws.Range("K7").Value = "0.22873154717987684" -as [double]
$ws.Range("I8").Value = "2.8735791784497015e-18" -as [double]
$ws.Range("J8").Value = "1.0503855952049763e-15" -as [double]
$ws.Range("K8").Value = "0.055435736427661746" -as [double]
$ws.Range("I9").Value = "-3.005364372432867e-18" -as [double]
$ws.Range("J9").Value = "4.9490305997717147e-17" -as [double]
$ws.Range("K9").Value = "0.14457374507012158" -as [double]
$ws.Range("I10").Value = "4.069423643307849e-18" -as [double]
$ws.Range("J10").Value = "4.654309318698814e-17" -as [double]
$ws.Range("K10").Value = "0.21872499874809398" -as [double]
$ws.Range("I11").Value = "-1.8983290101529263e-18" -as [double]
$ws.Range("J11").Value = "5.169079877678032e-17" -as [double]
$ws.Range("K11").Value = "0.22348287562520747" -as [double]
$ws.Range("G12").Value = "5.735429492448318e-17" -as [double]
$ws.Range("I12").Value = "6.92170227223703e-18" -as [double]
$ws.Range("J12").Value = "3.5286240803344516e-17" -as [double]
$ws.Range("K12").Value = "0.30937084288145605" -as [double]
$ws.Range("I13").Value = "1.1831871914550464e-18" -as [double]
$ws.Range("J13").Value = "4.507802533421588e-17" -as [double]
$ws.Range("K13").Value = "0.2572842811536125" -as [double]
$ws.Range("I14").Value = "3.0687940935842024e-17" -as [double]
$ws.Range("J14").Value = "6.210145228370588e-17" -as [double]
$ws.Range("K14").Value = "0.20858878311784557" -as [double]
$ws.Range("I15").Value = "-2.1000606693058588e-18" -as [double]
$ws.Range("J15").Value = "3.7832322617375345e-17" -as [double]
$ws.Range("K15").Value = "0.2586029544355496" -as [double]
$ws.Range("I16").Value = "-2.86470307352281e-17" -as [double]
$ws.Range("J16").Value = "6.596044955977111e-17" -as [double]
$ws.Range("K16").Value = "0.2775446012810049" -as [double]
$ws.Range("I17").Value = "1.123376921625402e-17" -as [double]
$ws.Range("J17").Value = "7.631314405505938e-17" -as [double]
$ws.Range("K17").Value = "0.22912256475322426" -as [double]
$ws.Range("I18").Value = "-2.8271233031310203e-18" -as [double]
$ws.Range("J18").Value = "1.1263454563329659e-16" -as [double]
$ws.Range("K18").Value = "0.2003993670096433" -as [double]
$ws.Range("I19").Value = "-5.398491720815776e-17" -as [double]
$ws.Range("J19").Value = "1.961615371560519e-16" -as [double]
$ws.Range("K19").Value = "0.24076528165076003" -as [double]
$ws.Range("I20").Value = "6.420313456529308e-18" -as [double]
$ws.Range("J20").Value = "1.1583416072701452e-16" -as [double]
$ws.Range("K20").Value = "0.22040604297113614" -as [double]
$ws.Range("I21").Value = "-1.6077639216319052e-17" -as [double]
$ws.Range("J21").Value = "1.0154558202708016e-16" -as [double]
$ws.Range("K21").Value = "0.19255606121190666" -as [double]
$ws.Range("I22").Value = "1.2078915123852518e-17" -as [double]
$ws.Range("J22").Value = "1.0499522712204066e-16" -as [double]
$ws.Range("K22").Value = "0.16469558534704862" -as [double]
$ws.Range("I23").Value = "-9.73372617075875e-18" -as [double]
$ws.Range("J23").Value = "6.661111989600624e-17" -as [double]
$ws.Range("K23").Value = "0.2430297233356031" -as [double]
$ws.Range("I24").Value = "5.0240383013267886e-18" -as [double]
$ws.Range("J24").Value = "1.1356664092924107e-16" -as [double]
$ws.Range("K24").Value = "0.16738053375939424" -as [double]
$ws.Range("I25").Value = "-1.1525426503307944e-17" -as [double]
$ws.Range("J25").Value = "1.1098160901425095e-16" -as [double]
$ws.Range("K25").Value = "0.18661120609878318" -as [double]
$ws.Range("I26").Value = "-2.2085599809889936e-18" -as [double]
$ws.Range("J26").Value = "5.1566091665115076e-17" -as [double]
$ws.Range("K26").Value = "0.23406163775290972" -as [double]
$ws.Range("I27").Value = "-2.994004221500712e-17" -as [double]
$ws.Range("J27").Value = "9.144249192873183e-17" -as [double]
$ws.Range("K27").Value = "0.2558883868433809" -as [double]
